$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Historical Times" sheet: add the new benchmark run (row 6) which
#    records the date of the run and the grand total time.
# ---------------------------------------------------------------------
$wsHist = $wb.Worksheets.Item("Historical Times")

# Give the new date cell the same short-date format as the cells above it.
$wsHist.Range("A6").NumberFormat = "m/d/yy"
$wsHist.Range("A6").Value = 40382
$wsHist.Range("B6").Value = 24471

# Move the selection down one row, as happens after entering data and
# pressing Enter.
$wsHist.Range("B7").Select()

# ---------------------------------------------------------------------
# 2. "Breakdown" sheet: update with the figures for the new run, and
#    mark it as the active sheet (the new run's breakdown is what the
#    user is now looking at).
# ---------------------------------------------------------------------
$wsBreak = $wb.Worksheets.Item("Breakdown")

$wsBreak.Range("B3").Value = 40382

$data = @(
    @(6,  13.5, 32,  0, 199.1,               244.5),
    @(7,  0.7,  0.4, 0, 310.2,               311.3),
    @(8,  3.7,  15.8,0, 500.3,               519.79999999999995),
    @(9,  0.3,  1.2, 0, 1003.7,              1005.2),
    @(10, 0.6,  2,   0, 930.9,               933.5),
    @(11, 1.1000000000000001, 1.9, 0, 333.8, 336.9),
    @(12, 0.2,  0.4, 0, 322.60000000000002,  323.2),
    @(13, 0.2,  0.2, 0, 222.1,               222.5),
    @(14, 0.2,  0.4, 0, 149.80000000000001,  150.4),
    @(15, 0.1,  0.1, 0, 283.5,               283.60000000000002),
    @(16, 0.2,  0.5, 0, 273.89999999999998,  274.7),
    @(17, 0.3,  0.4, 0, 153,                 153.6),
    @(18, 4.9000000000000004, 7.5, 0, 663.6, 676),
    @(19, 2.8,  3.8, 0, 505.1,               511.7),
    @(20, 1.2,  2.2999999999999998, 0, 313.3,316.7),
    @(21, 1.8,  2.4, 0, 427.5,               431.7),
    @(22, 6.1,  3.3, 0, 321.10000000000002,  330.6),
    @(23, 0.4,  0.8, 0, 317.8,               319),
    @(24, 0.3,  0.6, 0, 181.7,               182.6),
    @(25, 0.4,  0.7, 0, 252.1,               253.2),
    @(26, 1.6,  0.8, 0, 768.4,               770.8),
    @(27, 1.2,  1.1000000000000001, 0, 1222.5999999999999, 1224.8),
    @(28, 0.5,  1,   0, 484.7,               486.3),
    @(29, 4,    2.1, 0, 3284.7,              3290.8),
    @(30, 16,   3.9, 0, 3670.4,              3690.4),
    @(31, 0.5,  0.8, 0, 7226.7,              7228)
)

foreach ($row in $data) {
    $r = $row[0]
    $wsBreak.Cells.Item($r, 2).Value = $row[1]
    $wsBreak.Cells.Item($r, 3).Value = $row[2]
    $wsBreak.Cells.Item($r, 4).Value = $row[3]
    $wsBreak.Cells.Item($r, 5).Value = $row[4]
    $wsBreak.Cells.Item($r, 6).Value = $row[5]
}

# Apply a one-decimal thousands number format to the whole data block,
# matching the existing borders/fills already on each cell.
$wsBreak.Range("B6:E31").NumberFormat = "#,##0.0"
$wsBreak.Range("F6:F31").NumberFormat = "#,##0.0"
$wsBreak.Range("B32:F32").NumberFormat = "#,##0.0"

$wsBreak.PageSetup.Orientation = 1

$wsBreak.Activate()
$wsBreak.Range("B4").Select()
